# Add TP/FP/TN/FN classification columns (I:L) and summary counts (N:Q, S:T)
# to Sheet1 of the Vader bias-dataset workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "TP"
$ws.Range("J1").Value = "FP"
$ws.Range("K1").Value = "TN"
$ws.Range("L1").Value = "FN"

$ws.Range("N1").Value = "TP"
$ws.Range("O1").Value = "FP"
$ws.Range("P1").Value = "TN"
$ws.Range("Q1").Value = "FN"

# --- Row 2: first classification formulas, typed directly (not filled) ---
$ws.Range("I2").Formula = '=IF(AND(A2=-1,G2=-1),"TP")'
$ws.Range("J2").Formula = '=IF(AND(A2=1,G2=-1),"FP")'
$ws.Range("K2").Formula = '=IF(AND(A2=1,G2=1),"TN")'
$ws.Range("L2").Formula = '=IF(AND(A2=-1,G2=1),"FN")'

# --- Rows 3:66: fill the classification formulas down (creates one shared block) ---
$ws.Range("I3:I66").Formula = '=IF(AND(A3=-1,G3=-1),"TP")'
$ws.Range("J3:J66").Formula = '=IF(AND(A3=1,G3=-1),"FP")'
$ws.Range("K3:K66").Formula = '=IF(AND(A3=1,G3=1),"TN")'
$ws.Range("L3:L66").Formula = '=IF(AND(A3=-1,G3=1),"FN")'

# --- Rows 67:68: filled down separately afterwards (creates a second shared block) ---
$ws.Range("I67:I68").Formula = '=IF(AND(A67=-1,G67=-1),"TP")'
$ws.Range("J67:J68").Formula = '=IF(AND(A67=1,G67=-1),"FP")'
$ws.Range("K67:K68").Formula = '=IF(AND(A67=1,G67=1),"TN")'
$ws.Range("L67:L68").Formula = '=IF(AND(A67=-1,G67=1),"FN")'

# --- Summary counts in row 2 ---
$ws.Range("N2").Formula = '=COUNTIF(I2:I200,"TP")'
$ws.Range("O2").Formula = '=COUNTIF(J2:J200,"FP")'
$ws.Range("P2").Formula = '=COUNTIF(K2:K200,"TN")'
$ws.Range("Q2").Formula = '=COUNTIF(L2:L200,"FN")'

$ws.Range("S2").Formula = '=COUNTIF(A2:A100,1)'
$ws.Range("T2").Formula = '=COUNTIF(A2:A100,-1)'

# --- Final selection, matching the saved view state ---
$ws.Range("T2").Select()
